$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '51.544.27'
$ws.Range('E2').Value = '  +1.07%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.985.69'
$ws.Range('E3').Value = '  +1.46%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.12%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '382.22'
$ws.Range('E5').Value = '  +1.85%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '104.36'
$ws.Range('E6').Value = '  +3.23%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.546'
$ws.Range('E7').Value = '  +1.16%  '

$ws.Range('E8').Value = '  -0.03%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.597'
$ws.Range('E9').Value = '  +1.43%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.70'
$ws.Range('E10').Value = '  +1.11%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.138'
$ws.Range('E11').Value = '  -0.63%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0859'
$ws.Range('E12').Value = '  +1.11%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.456.26'
$ws.Range('E13').Value = '  +1.52%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.86'
$ws.Range('E14').Value = '  +3.18%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '18.53'
$ws.Range('E15').Value = '  +2.33%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.983.71'
$ws.Range('E16').Value = '  +1.48%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '11.26'
$ws.Range('E17').Value = '  +0.11%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.00'
$ws.Range('E18').Value = '  +0.43%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '51.541.85'
$ws.Range('E19').Value = '  +1.17%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '3.11'
$ws.Range('E20').Value = '  +0.96%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.61'
$ws.Range('E21').Value = '  +0.98%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.0₃0964'
$ws.Range('E22').Value = '  +0.85%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '70.41'
$ws.Range('E23').Value = '  +2.15%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '267.42'
$ws.Range('E24').Value = '  +0.39%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.23'
$ws.Range('E25').Value = '  +2.28%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.89'
$ws.Range('E26').Value = '  -3.10%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.28'
$ws.Range('E27').Value = '  -2.44%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.170'
$ws.Range('E28').Value = '  +4.00%  '

$ws.Range('E29').Value = '  +0.07%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '26.13'
$ws.Range('E30').Value = '  +1.74%  '

$ws.Range('E31').Value = '  -0.57%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '10.44'
$ws.Range('E32').Value = '  +4.40%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '34.74'
$ws.Range('E33').Value = '  +4.15%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '51.39'
$ws.Range('E34').Value = '  +0.83%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.07'
$ws.Range('E35').Value = '  +0.46%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0449'
$ws.Range('E36').Value = '  +1.43%  '

$ws.Range('E37').Value = '  -0.09%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.29'
$ws.Range('E38').Value = '  +4.87%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '16.97'
$ws.Range('E39').Value = '  +2.93%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.58'
$ws.Range('E40').Value = '  +4.64%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.117'
$ws.Range('E41').Value = '  +1.14%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.85'
$ws.Range('E42').Value = '  +2.17%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.85'
$ws.Range('E43').Value = '  +13.44%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '126.42'
$ws.Range('E44').Value = '  +3.26%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '21.45'
$ws.Range('E45').Value = '  +1.16%  '

$ws.Range('E46').Value = '  -0.03%  '

$ws.Range('B47').Value = 'TheGraph'
$ws.Range('C47').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.273'
$ws.Range('E47').Value = '  +0.80%  '

$ws.Range('B48').Value = 'ApeXProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.36'
$ws.Range('E48').Value = '  +1.05%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.027.71'
$ws.Range('E49').Value = '  +1.77%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.277.64'
$ws.Range('E50').Value = '  +1.17%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0329'
$ws.Range('E51').Value = '  +0.62%  '
